$d = $word.ActiveDocument

$replacements = @(
    @("49×84=", "29×42="),
    @("29×62=", "92×34="),
    @("50×26=", "63×14="),
    @("99×11=", "92×86="),
    @("60×93=", "43×89="),
    @("24×44=", "77×30="),
    @("29×97=", "18×75="),
    @("24×74=", "12×69="),
    @("43×38=", "52×24="),
    @("29×20=", "96×50="),
    @("15×55=", "79×88="),
    @("93×57=", "91×74="),
    @("97×18=", "34×98="),
    @("14×94=", "71×24="),
    @("68×46=", "75×12="),
    @("54×96=", "86×84="),
    @("70×48=", "77×64="),
    @("16×86=", "46×77="),
    @("52×87=", "80×16="),
    @("91×60=", "46×39="),
    @("30×33=", "99×53="),
    @("22×80=", "87×58="),
    @("43×43=", "31×79="),
    @("81×61=", "22×85="),
    @("98×48=", "87×31=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
